$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update random_forest row (row 2) metrics
$ws.Range("B2").Value = 3.6568046469664277
$ws.Range("C2").Value = 0.22932401822887949
$ws.Range("D2").Value = 2.7551307925967832
$ws.Range("E2").Value = 0.38285625817950447
$ws.Range("F2").Value = 0.61875379447685364
$ws.Range("G2").Value = 0.57626093358042985
$ws.Range("H2").Value = 0.61714374182049547
$ws.Range("I2").Value = 0.78566140842942933

# Update lsboost row (row 3) metrics
$ws.Range("B3").Value = 3.6487532235597304
$ws.Range("C3").Value = 0.22881910069941352
$ws.Range("D3").Value = 2.7624820257014844
$ws.Range("E3").Value = 0.38117219538465258
$ws.Range("F3").Value = 0.61739144421076375
$ws.Range("G3").Value = 0.57779851156520856
$ws.Range("H3").Value = 0.61882780461534748
$ws.Range("I3").Value = 0.78704235753092244

# Update neural_network row (row 4) metrics
$ws.Range("B4").Value = 3.8032154909667311
$ws.Range("C4").Value = 0.23850567442872173
$ws.Range("D4").Value = 2.8551359709596573
$ws.Range("E4").Value = 0.41412752492609234
$ws.Range("F4").Value = 0.64352740806129805
$ws.Range("G4").Value = 0.59717793599683799
$ws.Range("H4").Value = 0.58587247507390772
$ws.Range("I4").Value = 0.76549544100335287

# Row 5 (old_model) values remain unchanged
